$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "#AWS options" section (rows 27-34), written cell-by-cell in the
# same order the original author entered them so that the shared-string
# table comes out in the same sequence.
$ws.Range("A27").Value = "#AWS options"

$ws.Range("A28").Value = "AWS access key (confidential)"
$ws.Range("A29").Value = "AWS secret key (very confidential)"

$ws.Range("B28").Value = "accessKey"
$ws.Range("B29").Value = "secretKey"

$ws.Range("B30").Value = "region"
$ws.Range("A30").Value = "AWS region"

$ws.Range("A31").Value = "AWS AMI to be used"
$ws.Range("C31").Value = "ami-9dc5a68a"
$ws.Range("B31").Value = "AMI"

$ws.Range("A32").Value = "AWS instance type"
$ws.Range("B32").Value = "instanceType"
$ws.Range("C32").Value = "subnet-8ee22c7"

$ws.Range("B33").Value = "bootStorageSize"
$ws.Range("A33").Value = "How big boot image should be"
$ws.Range("C33").Value = "20GB"

$ws.Range("A34").Value = "AWS maximum number of instances to be used"
$ws.Range("B34").Value = "maxInstances"

# C34 keeps "5" as text (matches the other Text-formatted value cells,
# e.g. C9:C18), so force the Text number format before writing it.
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "5"

# Leave the selection where the author ended up.
$ws.Range("C35").Select()
